{"js": "// Renumber [[PERSON_N]] placeholders: every N >= 38 becomes N-1.\n// (The author removed one PERSON slot \u2014 formerly distinct PERSON_37/PERSON_38\n// referred to the same dative form, so PERSON_38 was dropped and every\n// higher-numbered placeholder shifted down by one.)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst THRESHOLD = 38;\nconst tokenRe = /\\[\\[PERSON_(\\d+)\\]\\]/g;\n\nfunction renumber(text) {\n  return text.replace(tokenRe, (match, numStr) => {\n    const n = parseInt(numStr, 10);\n    const newN = n >= THRESHOLD ? n - 1 : n;\n    return \"[[PERSON_\" + newN + \"]]\";\n  });\n}\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const para = items[i];\n  const oldText = para.text;\n  if (!oldText || oldText.indexOf(\"[[PERSON_\") === -1) continue;\n  const newText = renumber(oldText);\n  if (newText !== oldText) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Renumber [[PERSON_N]] placeholders: every N >= 38 becomes N-1.\n# (Formerly distinct PERSON_37/PERSON_38 referred to the same dative form,\n# so PERSON_38 was dropped and every higher-numbered placeholder shifted\n# down by one.)\n\n$d = $word.ActiveDocument\n\n$threshold = 38\n$pattern = [regex]'\\[\\[PERSON_(\\d+)\\]\\]'\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $r = $p.Range\n    $oldText = $r.Text\n    if ($oldText -eq $null) { continue }\n    # Strip the trailing paragraph mark before matching/rewriting.\n    $body = $oldText.TrimEnd([char]13, [char]7)\n    if ($body.IndexOf('[[PERSON_') -lt 0) { continue }\n\n    $matches = $pattern.Matches($body)\n    $newBody = $body\n    for ($j = $matches.Count - 1; $j -ge 0; $j--) {\n        $m = $matches[$j]\n        $n = [int]$m.Groups[1].Value\n        if ($n -ge $threshold) { $n = $n - 1 }\n        $repl = \"[[PERSON_\" + $n + \"]]\"\n        $newBody = $newBody.Substring(0, $m.Index) + $repl + $newBody.Substring($m.Index + $m.Length)\n    }\n\n    if ($newBody -ne $body) {\n        $r.Text = $newBody\n    }\n}\n"}
